# feat: add 2022-Q3 data
$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Update the "总计" (totals) sheet: insert a new first data row for
#    2022-Q3 and shift the previous quarters down by one row.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$totalRows = @(
    @(0, "2022-Q3", 9,  6.51),
    @(1, "2022-Q2", 14, 8.06),
    @(2, "2022-Q1", 25, 12.2),
    @(3, "2021-Q4", 5,  5.96),
    @(4, "2021-Q3", 4,  2.98),
    @(5, "2021-Q2", 1,  0.02)
)

$r = 2
foreach ($row in $totalRows) {
    $total.Cells.Item($r, 1).Value = $row[0]
    $total.Cells.Item($r, 2).Value = $row[1]
    $total.Cells.Item($r, 3).Value = $row[2]
    $total.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# ------------------------------------------------------------------
# 2) Insert the new "2022-Q3" sheet right after "总计" (i.e. before the
#    sheet that is currently "2022-Q2").
# ------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item(2)
$q3 = $wb.Worksheets.Add($beforeSheet)
$q3.Name = "2022-Q3"

# Visual formatting matching the other quarter sheets: bold, centered,
# thin-bordered header row + first column.
$headerRange = $q3.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$colA = $q3.Range("A2:A10")
$colA.Font.Bold = $true
$colA.HorizontalAlignment = -4108
$colA.VerticalAlignment = -4160
$colA.Borders.LineStyle = 1

# Columns that hold numeric-looking text values must be forced to Text
# so leading zeros (fund codes) and fixed decimal formatting are
# preserved exactly as strings instead of being parsed as numbers.
$q3.Range("B2:B10").NumberFormat = "@"
$q3.Range("D2:G10").NumberFormat = "@"

# Header row.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q3.Cells.Item(1, $c + 2).Value = $headers[$c]
}

# Data rows.
$q3Rows = @(
    @(0, "270002", "广发稳健增长混合A",       "160.46", "38.78", "2.20", "3.5301", 6),
    @(1, "009951", "广发稳健回报混合A",       "52.44",  "46.49", "4.09", "2.1448", 1),
    @(2, "009952", "广发稳健回报混合C",       "9.48",   "46.49", "4.09", "0.3877", 1),
    @(3, "501070", "广发睿阳三年定期开放混合", "6.62",   "51.01", "4.58", "0.3032", 4),
    @(4, "008602", "方正富邦新兴成长混合A",   "1.23",   "86.03", "4.89", "0.0601", 2),
    @(5, "009326", "广发稳健增长混合C",       "2.30",   "38.78", "2.20", "0.0506", 6),
    @(6, "015032", "中融医药消费混合A",       "0.54",   "90.81", "4.95", "0.0267", 4),
    @(7, "008603", "方正富邦新兴成长混合C",   "0.03",   "86.03", "4.89", "0.0015", 2),
    @(8, "015033", "中融医药消费混合C",       "0.02",   "90.81", "4.95", "0.0010", 4)
)

$r = 2
foreach ($row in $q3Rows) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $q3.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}

Write-Output "2022-Q3 sheet added"
